$d = $word.ActiveDocument

# 1. Update the "Curso (semestre ideal)" line to add EA (1)
$d.Content.Find.Execute("Curso (semestre ideal): EP (10)", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Curso (semestre ideal): EA (1), EP (10)", 2)

# 2. Remove the trailing "Requisitos" heading paragraph and the LOB1009 bullet paragraph
$count = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($count)
$pHeading = $d.Paragraphs.Item($count - 1)
$r = $d.Range($pHeading.Range.Start, $pLast.Range.End)
$r.Delete()
